# Fixed anonymous sessions bug: add a new "Anonymous" / "Anônimo" text-db
# entry into the translation table on sheet1, right before the
# "TelaDeRegistro" section header (which currently lives at row 16).
#
# This inserts a new row 16 (pushing the header and everything below it
# down by one row) and populates it with the English key and the
# Portuguese translation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16, shifting existing rows 16-95 down
# to 17-96.
$ws.Rows("16:16").Insert() | Out-Null

# Populate the newly inserted row with the new translation pair.
$ws.Range("A16").Value = "Anonymous"
$ws.Range("B16").Value = "Anônimo"

# Match the saved selection state (Excel's last active cell moved to A16).
$ws.Range("A16").Select() | Out-Null
